$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.213.52'
$ws.Range('E2').Value = '  +5.34%  '
$ws.Range('D3').Value = '3.333.07'
$ws.Range('E3').Value = '  +2.11%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  -0.18%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '412.19'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +3.85%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '110.97'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +1.08%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.584'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  +3.85%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.999'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  -0.23%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.632'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +1.16%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '39.36'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -0.08%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0986'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +2.40%  '
$ws.Range('E12').Value = '  +1.07%  '
$ws.Range('D13').Value = '3.856.94'
$ws.Range('E13').Value = '  +1.88%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '8.41'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +2.71%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '19.64'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +2.84%  '
$ws.Range('D16').Value = '3.350.46'
$ws.Range('E16').Value = '  +2.53%  '
$ws.Range('E17').Value = '  +0.18%  '
$ws.Range('D18').Value = '59.919.77'
$ws.Range('E18').Value = '  +5.17%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '10.77'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +0.29%  '
$ws.Range('E20').Value = '  +2.10%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.0000110'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +3.92%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '13.22'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +2.41%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '299.71'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -1.94%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '75.05'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -0.22%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '3.17'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +0.46%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '8.12'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +3.06%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '28.56'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +1.49%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.78'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +6.96%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '4.47'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +1.73%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.179'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +5.57%  '
$ws.Range('E31').Value = '  +4.32%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.58'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +20.57%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '11.47'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +4.18%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.00'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +0.25%  '
$ws.Range('E35').Value = '  +4.55%  '
$ws.Range('E36').Value = '  +4.27%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '52.22'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +1.43%  '
$ws.Range('B38').Value = 'Stacks'
$ws.Range('C38').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.10'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -0.64%  '
$ws.Range('B39').Value = 'FirstDigitalUSD'
$ws.Range('C39').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.998'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +0.07%  '
$ws.Range('E40').Value = '  -3.56%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '137.96'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +2.81%  '
$ws.Range('E42').Value = '  +2.41%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.294'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +4.71%  '
$ws.Range('E44').Value = '  -0.38%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.93'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -1.41%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '16.87'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -1.88%  '
$ws.Range('E47').Value = '  +8.44%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '22.32'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +0.91%  '
$ws.Range('D49').Value = '2.183.34'
$ws.Range('E49').Value = '  +1.56%  '
$ws.Range('E50').Value = '  +1.63%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.01'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -0.70%  '
